# [Improvement] On terminology : room -> bed
#
# Renames the "rooms" worksheet to "beds" and updates the room/bed
# terminology used in its header row, then makes it the active sheet
# (matching the new selection left behind in the source workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rooms")

# Header row terminology: rooms -> beds
$ws.Range("A1").Value = "all_beds"
$ws.Range("B1").Value = "new_beds"
$ws.Range("C1").Value = "old_beds"
$ws.Range("E1").Value = "new_beds_service"
$ws.Range("F1").Value = "old_beds_service"
$ws.Range("G1").Value = "beds_capacities"

# Rename the sheet itself
$ws.Name = "beds"

# This sheet becomes the active tab/selection (previously "babies" was active)
[void]$ws.Activate()
[void]$ws.Range("N20").Select()
